# Update "想去人数" (want-to-go headcount) figures in the "展览" sheet
# and the mirrored "全部类型" sheet, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value  = 221
$wsExhibit.Range("F3").Value  = 525
$wsExhibit.Range("F4").Value  = 121
$wsExhibit.Range("F6").Value  = 55
$wsExhibit.Range("F7").Value  = 73
$wsExhibit.Range("F8").Value  = 7162
$wsExhibit.Range("F9").Value  = 261
$wsExhibit.Range("F10").Value = 411
$wsExhibit.Range("F11").Value = 3630
$wsExhibit.Range("F13").Value = 524
$wsExhibit.Range("F16").Value = 85

# Sheet 4: 全部类型 (All types) — same events, different row offsets
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F3").Value  = 221
$wsAll.Range("F4").Value  = 525
$wsAll.Range("F5").Value  = 121
$wsAll.Range("F7").Value  = 55
$wsAll.Range("F8").Value  = 73
$wsAll.Range("F10").Value = 7162
$wsAll.Range("F12").Value = 261
$wsAll.Range("F13").Value = 411
$wsAll.Range("F14").Value = 3630
$wsAll.Range("F16").Value = 524
$wsAll.Range("F19").Value = 85
